# Rename the inline picture shapes' display Name (the wp:docPr / pic:cNvPr
# "name" attribute) for the three logo images touched by this revision:
#   - the Pearson logo in both footers: image1.png -> image2.png
#   - the BTEC logo in the header:      image2.jpg -> image1.jpg
#
# Word's InlineShape object has no writable .Name property (only Shape /
# ShapeRange expose Name), so we convert each inline picture to a floating
# Shape, rename it, then convert it straight back to an inline picture -
# this round trip is what actually updates the drawing's docPr/name in the
# underlying OOXML while keeping the shape inline (wp:inline, not anchored).

$d = $word.ActiveDocument

function Rename-InlineLogo($range, $targetName) {
    if ($range.InlineShapes.Count -gt 0) {
        $shp = $range.InlineShapes.Item(1)
        $floating = $shp.ConvertToShape()
        $floating.Name = $targetName
        [void]$floating.ConvertToInlineShape()
    }
}

foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            Rename-InlineLogo $hdr.Range "image1.jpg"
        }
    }
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            Rename-InlineLogo $ftr.Range "image2.png"
        }
    }
}
